$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> (new Iteration Path value in column J, new Target Date serial value in column P)
$changes = @{
    8   = @{ J = "Hydrangea";      P = 45374 }
    11  = @{ J = "Marigold";       P = 45703 }
    48  = @{ J = "Rose";           P = 45661 }
    62  = @{ J = "Zinnia";         P = 45458 }
    84  = @{ J = "Lily";           P = 45682 }
    99  = @{ J = "Anemone";        P = 45777 }
    103 = @{ J = "Chrysanthemum";  P = 45521 }
    105 = @{ J = "Hydrangea";      P = 45374 }
    115 = @{ J = "Dahlia";         P = 45395 }
    118 = @{ J = "Chrysanthemum";  P = 45521 }
    119 = @{ J = "Daffodil";       P = 45479 }
    121 = @{ J = "Jasmine";        P = 45437 }
    132 = @{ J = "Sunflower";      P = 45416 }
    134 = @{ J = "Orchid";         P = 45584 }
    136 = @{ J = "Orchid";         P = 45584 }
    137 = @{ J = "Chrysanthemum";  P = 45521 }
    139 = @{ J = "Lavender";       P = 45500 }
    146 = @{ J = "Orchid";         P = 45584 }
    148 = @{ J = "Hibiscus";       P = 45542 }
    150 = @{ J = "Orchid";         P = 45584 }
    152 = @{ J = "Zinnia";         P = 45458 }
    153 = @{ J = "Jasmine";        P = 45437 }
}

foreach ($row in $changes.Keys) {
    $entry = $changes[$row]
    $ws.Range("J$row").Value = $entry.J
    $ws.Range("P$row").Value = $entry.P
}
